$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 57, pushing the old row 57 (and below) down to 58.
# The existing row 56 stays put for now.
$ws.Rows.Item(57).Insert()

# Duplicate the (still unmodified) row 56 into the newly-created row 57, so
# row 57 now holds the same "Especial" record that used to live at row 56.
$ws.Range("A56:T56").Copy()
$ws.Range("A57").PasteSpecial()

# Now overwrite row 56 with the updated "Especial" record values.
$ws.Range("D56").Value = 44516
$ws.Range("M56").Value = 250
$ws.Range("N56").Value = 25000
$ws.Range("O56").Value = 25000
$ws.Range("P56").Value = 25000
$ws.Range("S56").Value = 2500
